$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column K: "Authorship Resource" header + value for every data row
$ws.Range("K1").Value = "Authorship Resource"

$authors = "Daniela Subotic, Noémi Villars-Amberg"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 11).Value = $authors
}

# Set a reasonable width for the new column (matches author's bestFit width as closely as possible)
$ws.Columns.Item(11).ColumnWidth = 31.8

# Update selection to reflect the newly filled column, matching the authored edit
$ws.Range("K2:K11").Select() | Out-Null
